$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19: sd / Statement-non-opinion -> b / Acknowledge (Backchannel)
$ws.Range("I19").Value = "b"
$ws.Range("J19").Value = "Acknowledge (Backchannel)"

# Row 20: % / Uninterpretable -> sd / Statement-non-opinion
$ws.Range("I20").Value = "sd"
$ws.Range("J20").Value = "Statement-non-opinion"

# Row 64: b / Acknowledge (Backchannel) -> aa / Agree/Accept
$ws.Range("I64").Value = "aa"
$ws.Range("J64").Value = "Agree/Accept"

# Row 72: sv / Statement-opinion -> sd / Statement-non-opinion
$ws.Range("I72").Value = "sd"
$ws.Range("J72").Value = "Statement-non-opinion"

# Row 73: sv / Statement-opinion -> sd / Statement-non-opinion
$ws.Range("I73").Value = "sd"
$ws.Range("J73").Value = "Statement-non-opinion"
